# Updated code quality rules and rel rating
#
# The "Rules" sheet lists CodeQuality rules. The rule formerly shown in row 35
# ("BannedPaths") is removed from that position (rows 36-40 shift up to
# 35-39), and a renamed/updated version of that rule ("BannedPath", severity
# "Critical", no Tags) is (re)inserted at row 40 -- effectively moving/
# resorting that single rule further down the list and tweaking its data.
#
# All of the other cell-index churn visible in the raw XML diff (rows 71-74,
# 116-132) is a pure side-effect of one shared string ("BannedPaths") having
# been removed from the middle of the sharedStrings table -- the visible text
# in those rows is unchanged, so no action is needed for them here; Excel's
# own shared-string de-duplication on save reproduces that renumbering for us
# automatically once we edit the text of the one cell that really changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row (currently row 35: "BannedPaths" / .../libs / Bug / Blocker)
# Rows 36..134 shift up by one (row 36 -> 35, ..., row 41 -> 40, ...).
$ws.Rows("35").Delete()

# Make room for the updated rule at its new location, row 40, pushing the
# row that is currently there (and everything below) back down by one.
$ws.Rows("40").Insert()

# Populate the re-inserted / renamed rule row.
$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"
# (No Tags value for this row -- column E is left blank.)

# Update the sheet's saved selection from B38 to A37.
$ws.Range("A37").Select() | Out-Null
